# Add two new columns, I ("I0") and J ("IF"), to the worksheet.
# I0 is always 1; IF mirrors the existing IP (column H) value for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered alignment) from
# the existing H1 header cell onto the two new header cells so the style
# matches the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2 .. 42) ---------------------------------------------------
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2   # column H ("IP") value
    $ws.Cells.Item($r, 9).Value = 1            # column I ("I0")
    $ws.Cells.Item($r, 10).Value = $ipValue    # column J ("IF")
}
